$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "OSRM Raw" column (M) -------------------------------------------
$ws.Range("M1").Value = "OSRM Raw"
$ws.Range("M1").Font.Bold = $true

$ws.Range("M2").Value = 77.400000000000006
$ws.Range("M3").Value = 78.5
$ws.Range("M4").Value = 86.2
$ws.Range("M5").Value = 79.2
$ws.Range("M6").Value = 79.5
$ws.Range("M7").Formula = "=AVERAGE(M2:M6)"

# --- Fix minor numbering error (no net routing change) --------------------
$ws.Range("C35").Value = 17
$ws.Range("C38").Value = 27

# --- Drop now-unused, style-only P:Q columns so the used range shrinks ----
$ws.Range("P1:Q49").Clear()

# --- Selection / view state -------------------------------------------
$ws.Range("M7").Select() | Out-Null
